# Auto-generated files on 2026-01-14
# Update the Hot Stock Top20 table contents (rows 2-21, columns A-C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("海格通信", "岩山科技", "华胜天成"),
    @("利欧股份", "利欧股份", "岩山科技"),
    @("岩山科技", "特变电工", "海格通信"),
    @("特变电工", "山子高科", "利欧股份"),
    @("三维通信", "海格通信", "航天发展"),
    @("蓝色光标", "航天发展", "金风科技"),
    @("金风科技", "金风科技", "山子高科"),
    @("美年健康", "蓝色光标", "神剑股份"),
    @("山子高科", "航天电子", "特变电工"),
    @("航天发展", "华胜天成", "雷科防务"),
    @("华胜天成", "三维通信", "省广集团"),
    @("雷科防务", "雷科防务", "蓝色光标"),
    @("省广集团", "中国卫星", "天下秀"),
    @("神剑股份", "省广集团", "国晟科技"),
    @("航天电子", "美年健康", "三维通信"),
    @("天下秀",   "东方财富", "中国卫通"),
    @("中国卫通", "神剑股份", "美年健康"),
    @("值得买",   "易点天下", "银河电子"),
    @("中国卫星", "中国卫通", "中国一重"),
    @("石基信息", "天下秀",   "航天电子")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
